$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update jurusan (column E) for row 4 from "TKJ" to new value "tbsm"
$ws.Range("E4").Value = "tbsm"

# Update d_kelas (column F) for row 5 from 1 to 4
$ws.Range("F5").Value = 4

# Move the active selection to F5
$ws.Range("F5").Select()
